$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.291.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5291"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2656"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07840"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.533"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.664.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.898.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8110"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.304.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.733"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.92%  "

$ws.Range("E21").Value = "  +4.24%  "

$ws.Range("E22").Value = "  +1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.075"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.008"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1220"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.529"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05913"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.283"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.520"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.337"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.46%  "

$ws.Range("E34").Value = "  -0.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9647"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.69%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5811"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01613"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.988"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.077.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8576"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.807.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.014"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4410"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.094"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₈102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05144"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
